# Auto-generated edits for Leviathan_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N)
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 371.66666
$ws.Range("I2").Value = 589.2
$ws.Range("J2").Value = 99.75
$ws.Range("K2").Value = 589.2
$ws.Range("L2").Value = 99.75
$ws.Range("M2").Value = -476.2
$ws.Range("N2").Value = -325.75

# Row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 11424.833
$ws.Range("I43").Value = 11750
$ws.Range("J43").Value = 11262.25
$ws.Range("K43").Value = 11750
$ws.Range("L43").Value = 11262.25
$ws.Range("M43").Value = -11681
$ws.Range("N43").Value = -11400.25

# Row 62 (Leve Item ID 27781)
$ws.Range("H62").Value = 72039
$ws.Range("I62").Value = 88048.25
$ws.Range("J62").Value = 8002
$ws.Range("K62").Value = 88048.25
$ws.Range("L62").Value = 8002
$ws.Range("M62").Value = -87424.25
$ws.Range("N62").Value = -9250

# Row 65 (Leve Item ID 27781)
$ws.Range("H65").Value = 72039
$ws.Range("I65").Value = 88048.25
$ws.Range("J65").Value = 8002
$ws.Range("K65").Value = 440241.25
$ws.Range("L65").Value = 40010
$ws.Range("M65").Value = -437121.25
$ws.Range("N65").Value = -46250

# Row 135 (Leve Item ID 44047)
$ws.Range("H135").Value = 34344.266
$ws.Range("I135").Value = 924.8
$ws.Range("J135").Value = 101183.2
$ws.Range("K135").Value = 8323.199999999999
$ws.Range("L135").Value = 910648.7999999999
$ws.Range("M135").Value = -5788.199999999999
$ws.Range("N135").Value = -915718.7999999999


$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 52899.66
$ws.Range("I32").Value = 30030.334
$ws.Range("J32").Value = 217558.8
$ws.Range("K32").Value = 30030.334
$ws.Range("L32").Value = 217558.8
$ws.Range("M32").Value = -29743.334
$ws.Range("N32").Value = -218132.8

# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 773045.9399999999
$ws.Range("I45").Value = 2004301.4
$ws.Range("K45").Value = 2004301.4
$ws.Range("M45").Value = -2003924.4

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 1887.65
$ws.Range("I132").Value = 1316.2
$ws.Range("K132").Value = 3948.6
$ws.Range("M132").Value = -1418.6


$ws = $wb.Worksheets.Item("BSM")
# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 2101.7932
$ws.Range("I86").Value = 1884.5555
$ws.Range("K86").Value = 1884.5555
$ws.Range("M86").Value = -761.5554999999999

# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 2101.7932
$ws.Range("I89").Value = 1884.5555
$ws.Range("K89").Value = 9422.7775
$ws.Range("M89").Value = -3806.7775

# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 1322.4166
$ws.Range("I99").Value = 644.8889
$ws.Range("K99").Value = 644.8889
$ws.Range("M99").Value = 853.1111

# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 3575606.5
$ws.Range("I105").Value = 3850383
$ws.Range("K105").Value = 3850383
$ws.Range("M105").Value = -3848636

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 2296
$ws.Range("I134").Value = 1701.7693
$ws.Range("K134").Value = 5105.3079
$ws.Range("M134").Value = -2570.3079


$ws = $wb.Worksheets.Item("CRP")
# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 3286.7144
$ws.Range("I16").Value = 4527.5
$ws.Range("K16").Value = 4527.5
$ws.Range("M16").Value = -4240.5

# Row 22 (Leve Item ID 5367)
$ws.Range("H22").Value = 240
$ws.Range("I22").Value = 240
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 240
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 110
$ws.Range("N22").ClearContents()

# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 2473.3333
$ws.Range("I31").Value = 2473.3333
$ws.Range("K31").Value = 2473.3333
$ws.Range("M31").Value = -2178.3333

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2473.3333
$ws.Range("I34").Value = 2473.3333
$ws.Range("K34").Value = 2473.3333
$ws.Range("M34").Value = -2271.3333

# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 1983.3334
$ws.Range("I58").Value = 1975
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 1975
$ws.Range("L58").Value = 2000
$ws.Range("M58").Value = -1772
$ws.Range("N58").Value = -2406

# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 3286.7144
$ws.Range("I113").Value = 4527.5
$ws.Range("K113").Value = 4527.5
$ws.Range("M113").Value = -2357.5

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2098.348
$ws.Range("I132").Value = 2156.2856
$ws.Range("K132").Value = 6468.8568
$ws.Range("M132").Value = -3938.8568

# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 1983.3334
$ws.Range("I136").Value = 1975
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 5925
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -3375
$ws.Range("N136").Value = -11100


$ws = $wb.Worksheets.Item("CUL")
# Row 11 (Leve Item ID 4745)
$ws.Range("H11").Value = 537
$ws.Range("I11").Value = 499.5
$ws.Range("K11").Value = 1498.5
$ws.Range("M11").Value = -1358.5

# Row 48 (Leve Item ID 4724)
$ws.Range("H48").Value = 2000
$ws.Range("I48").Value = 2000
$ws.Range("K48").Value = 6000
$ws.Range("M48").Value = -5750

# Row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 2707.25
$ws.Range("J68").Value = 2943
$ws.Range("L68").Value = 8829
$ws.Range("N68").Value = -10451

# Row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 2707.25
$ws.Range("J71").Value = 2943
$ws.Range("L71").Value = 26487
$ws.Range("N71").Value = -34599


$ws = $wb.Worksheets.Item("GSM")
# Row 62 (Leve Item ID 11983)
$ws.Range("H62").Value = 22085
$ws.Range("J62").Value = 22085
$ws.Range("L62").Value = 22085
$ws.Range("N62").Value = -23457

# Row 65 (Leve Item ID 11983)
$ws.Range("H65").Value = 22085
$ws.Range("J65").Value = 22085
$ws.Range("L65").Value = 66255
$ws.Range("N65").Value = -73119

# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 2494.3333
$ws.Range("I132").Value = 2569.8096
$ws.Range("J132").Value = 1966
$ws.Range("K132").Value = 7709.4288
$ws.Range("L132").Value = 5898
$ws.Range("M132").Value = -5179.4288
$ws.Range("N132").Value = -10958


$ws = $wb.Worksheets.Item("LTW")
# Row 2 (Leve Item ID 2631)
$ws.Range("H2").Value = 500
$ws.Range("I2").Value = 500
$ws.Range("K2").Value = 500
$ws.Range("M2").Value = -388

# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 29183
$ws.Range("I7").Value = 38976
$ws.Range("K7").Value = 38976
$ws.Range("M7").Value = -38864

# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 1024.75

# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 1024.75

# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 29183
$ws.Range("I126").Value = 38976
$ws.Range("K126").Value = 116928
$ws.Range("M126").Value = -114458

# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 2437.1353
$ws.Range("I136").Value = 1876.2084
$ws.Range("K136").Value = 5628.6252
$ws.Range("M136").Value = -3078.6252


$ws = $wb.Worksheets.Item("WVR")
# Row 12 (Leve Item ID 3316)
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
